# Applies the "Recycling Challenge UI" spec update:
#  1. Clarify scanning is once-per-station-per-day (several sentences).
#  2. Drop the stray _GoBack bookmark that used to sit after the
#     "...positive impact on the earth's environment." sentence.
#  3. Clarify the "already completed today" message wording.
#  4. Add a new bullet describing the backend tracking requirement.
#  5. Clarify the "only once a day" -> "once per recycling station per day" bullet.
#  6. Reword the Given/When example lines for the repeat-scan scenario, and
#     re-home the _GoBack bookmark at the end of the new "When:" sentence
#     (Word leaves it at the last edit point).

$d = $word.ActiveDocument

# --- 1. Scan-limit description under the rules bullet ---------------------
$d.Content.Find.Execute(
    "The user can complete the scan at the most once a day. And the repeated scanning is invalid.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The user can complete the scan once per recycling station per day. Repeated scans at the same station on the same day are invalid.",
    1) | Out-Null

# --- 2. Remove the old _GoBack bookmark (it gets re-added later) ----------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 3. "already completed" popup trigger wording ---------------------------
# (there are two near-identical sentences in the doc; only the first - the
# one that is not preceded by "is displayed. " - should change)
$d.Content.Find.Execute(
    "If the user scans the code but has completed the task today, the system will display ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If the user scans the code from a station they have already scanned today, the system will display ",
    1) | Out-Null

# --- 4. New bullet about backend tracking, inserted right before the
#        "Result calculation" bullet -----------------------------------------
$rng = $d.Content
$rng.Find.Execute("Result calculation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$resultPara = $rng.Paragraphs(1)
$insertPoint = $resultPara.Range
$insertPoint.Collapse(1)
$insertPoint.InsertBefore("The backend must track which station codes the user has scanned each day to ensure the user can only scan once per station per day.`r")

$rng2 = $d.Content
$rng2.Find.Execute("The backend must track which station codes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$newPara = $rng2.Paragraphs(1)
$newPara.Range.ListFormat.RemoveNumbers()
$newPara.Range.ParagraphFormat.LeftIndent = 36
$newPara.Range.ParagraphFormat.CharacterUnitLeftIndent = 0
$newPara.Range.ParagraphFormat.FirstLineIndent = 24
$newPara.Range.ParagraphFormat.CharacterUnitFirstLineIndent = 200

# --- 5. "only once a day" bullet -------------------------------------------
$d.Content.Find.Execute(
    "   The user can complete this task only once a day.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "   The user can complete this task once per recycling station per day.",
    1) | Out-Null

# --- 6. Given/When example text --------------------------------------------
$d.Content.Find.Execute(
    " The user has successfully completed this task once today.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " The user has successfully scanned at this specific recycling station today.",
    1) | Out-Null

# Replace the "When:" sentence, keeping the trailing period intact for now
# (splitting it off right at a paragraph end confuses the bookmark anchor),
# then move the _GoBack bookmark to sit right before that trailing period -
# this naturally splits "...same day." into "...same day" + "." runs with
# the bookmark in between, matching how Word re-homes _GoBack at the last
# edited spot.
$d.Content.Find.Execute(
    " The user scans code again.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " The user scans the QR code of the same station again on the same day.",
    1) | Out-Null

$tail = $d.Content
$tail.Find.Execute("on the same day.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$periodPos = $tail.End - 1
$word.Selection.SetRange($periodPos, $periodPos)
$d.Bookmarks.Add("_GoBack", $word.Selection.Range)

Write-Output "Recycling Challenge UI updated."
